# Daily attendance processing - 2026-02-01 15:39:14
# Swap the order of names in the "Recorded By" (column G) cells that
# currently read "dnasr281@gmail.com, System" so they read
# "System, dnasr281@gmail.com" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
}
